$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row of results for "DORN_nohints" (SGD method, no albedo/no falloff/no noise)
$ws.Range("A21").Value = "DORN_nohints"
$ws.Range("B21").Value = 0.97945953669757
$ws.Range("C21").Value = 0.994678041453117
$ws.Range("D21").Value = 0.998067319060799
$ws.Range("E21").Value = 0.118581589964839
$ws.Range("F21").Value = 0.309140993793327
$ws.Range("G21").Value = 0.0840228550213256
$ws.Range("H21").Value = 0.031475905720026
$ws.Range("I21").Value = 0.038626020417079
$ws.Range("J21").Value = 0.105912321823124

# Relabel last section's DORN row to DORN_nohints as well
$ws.Range("A31").Value = "DORN_nohints"

# Update the active selection to A32
$ws.Range("A32").Select()
